# draft-gandhi-ippm-stamp-direct-00.pptx - "Add files via upload"
#
# Slide 5 ("Direct Measurement TLV vs. Direct Measurement Test Packet")
# contains a comparison table. This edit:
#   1) Swaps the widths of the 2nd and 3rd columns.
#   2) Adds a follow-up line "(out-of-order data packet support)" under the
#      "Alternate-marking method packet loss - using block number for
#      counters" row.
#   3) Tightens up every row's height (a side-effect of PowerPoint's table
#      autofit recalculating row heights once the column widths/text
#      change).

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(5)

# Find the table shape on the slide.
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
        break
    }
}

$tbl = $tableShape.Table

# 1) Swap column 2 / column 3 widths (in points; stored EMU / 914400 * 72).
$tbl.Columns.Item(2).Width = 2438400 / 12700
$tbl.Columns.Item(3).Width = 2667000 / 12700

# 2) Append the extra line to the "Alternate-marking method..." row (row 9,
#    first column) as a new paragraph.
$cell = $tbl.Cell(9, 1)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Text = $tr.Text + "`r(out-of-order data packet support)"

# 3) Apply the recalculated row heights (points = EMU / 12700).
$rowHeights = @{
    1  = 241541
    2  = 357587
    3  = 241541
    4  = 498067
    5  = 357587
    6  = 357587
    7  = 241541
    8  = 241541
    9  = 498067
    10 = 241541
}
foreach ($rowIndex in $rowHeights.Keys) {
    $tbl.Rows.Item([int]$rowIndex).Height = $rowHeights[$rowIndex] / 12700
}
